# "Generate Report for Handoff" -- refresh the localization-status report
# with a new handoff cycle: new source file UUIDs, a new status text, new
# handoff timestamps, and collapsed "Latest Target File"/"Latest Handback
# File" columns (no handback has happened yet for this cycle).

$wb = $excel.ActiveWorkbook

$oldMd1 = "6d09ef12-6522-4ff2-bb86-40ce91ed510b.md"
$oldMd2 = "91bfec59-e124-4100-ba93-f5eb2b019972.md"
$newMd1 = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$newMd2 = "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md"

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: only the hyperlink display text needs to move to the
# new file names (cell values already reference the same shared-string
# slots, which we're updating in place).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$rA2 = $wsOverview.Range("A2")
$addrA2 = $rA2.Hyperlinks.Item(1).Address
$rA2.Hyperlinks.Delete()
$rA2.Hyperlinks.Add($rA2, $addrA2, "", "", $newMd1)

$rA3 = $wsOverview.Range("A3")
$addrA3 = $rA3.Hyperlinks.Item(1).Address
$rA3.Hyperlinks.Delete()
$rA3.Hyperlinks.Add($rA3, $addrA3, "", "", $newMd2)

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper applied identically to the "zh-cn" and "de-de" detail sheets.
# ---------------------------------------------------------------------
function Update-DetailSheet($sheetName, $newXlf, $handoffTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $oldXlf1 = $ws.Range("C2").Value
    $oldXlf2 = $ws.Range("C3").Value

    # Row 2 (was 6d09ef12...)
    $rA2 = $ws.Range("A2")
    $linkA2 = $rA2.Hyperlinks.Item(1).Address
    $rA2.Hyperlinks.Delete()
    $rA2.Hyperlinks.Add($rA2, $linkA2, "", "", $newMd1)

    $rC2 = $ws.Range("C2")
    $linkC2 = $rC2.Hyperlinks.Item(1).Address
    $rC2.Hyperlinks.Delete()
    $rC2.Value = $newXlf
    $rC2.Hyperlinks.Add($rC2, $linkC2, "", "", $newXlf)

    $ws.Range("B2").Value = $newStatus
    $ws.Range("D2").Value = $handoffTime
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Include"
    $ws.Range("E2:F2").Clear()
    $ws.Range("E2").Hyperlinks.Delete()
    $ws.Range("F2").Hyperlinks.Delete()

    # Row 3 (was 91bfec59...)
    $rA3 = $ws.Range("A3")
    $linkA3 = $rA3.Hyperlinks.Item(1).Address
    $rA3.Hyperlinks.Delete()
    $rA3.Hyperlinks.Add($rA3, $linkA3, "", "", $newMd2)

    $rC3 = $ws.Range("C3")
    $linkC3 = $rC3.Hyperlinks.Item(1).Address
    $rC3.Hyperlinks.Delete()
    $rC3.Value = $newXlf
    $rC3.Hyperlinks.Add($rC3, $linkC3, "", "", $newXlf)

    $ws.Range("B3").Value = $newStatus
    $ws.Range("D3").Value = $handoffTime
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"
    $ws.Range("E3:F3").Clear()
    $ws.Range("E3").Hyperlinks.Delete()
    $ws.Range("F3").Hyperlinks.Delete()

    # Row 4 (.localization-config) -- text unchanged, only the
    # shared-string indices shift because of the removals above.
    $ws.Range("D4").Value = "0001-01-01 00:00:00"
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Ignored"
}

Update-DetailSheet "zh-cn" "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf" "2016-03-10 19:06:11"
Update-DetailSheet "de-de" "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf" "2016-03-10 19:06:16"
